$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before old row 3, so old row3 -> row5, old row4 -> row6
$ws.Rows.Item(3).Resize(2).Insert()

# Row 2
$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', None),
                (''selector'',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                (''model'',
                 BaggingClassifier(estimator=LogisticRegression(C=0.001,
                                                                max_iter=1000,
                                                                random_state=42,
                                                                solver=''saga''),
                                   n_estimators=50, random_state=42))])'
$ws.Range("B2").Value = 0.6952380952380952
$ws.Range("C2").Value = '{''selector'': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), ''scaler'': None, ''model__n_estimators'': 50, ''model__estimator__solver'': ''saga'', ''model__estimator__penalty'': ''l2'', ''model__estimator__class_weight'': None, ''model__estimator__C'': 0.001}'
$ws.Range("D2").Value = 0.3333333333333333
$ws.Range("E2").Value = '[1 1 0 0 1 0 0 0 0 1 0 1]'
$ws.Range("F2").Value = '[0 1 1 0 0 1 1 0 1 0 1 1]'
$ws.Range("G2").Value = 77
$ws.Range("H2").Value = 0.7115238095238096
$ws.Range("I2").Value = 0.02434195695677412
$ws.Range("J2").Value = 0.5925714285714286
$ws.Range("K2").Value = 0.04926213077368703

# Row 3
$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', None),
                (''selector'',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                (''model'',
                 BaggingClassifier(estimator=LogisticRegression(C=3,
                                                                class_weight=''balanced'',
                                                                max_iter=1000,
                                                                random_state=42,
                                                                solver=''saga''),
                                   n_estimators=50, random_state=42))])'
$ws.Range("B3").Value = 0.6571428571428571
$ws.Range("C3").Value = '{''selector'': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), ''scaler'': None, ''model__n_estimators'': 50, ''model__estimator__solver'': ''saga'', ''model__estimator__penalty'': ''l2'', ''model__estimator__class_weight'': ''balanced'', ''model__estimator__C'': 3}'
$ws.Range("D3").Value = 0.9333333333333333
$ws.Range("E3").Value = '[1 1 0 1 0 0 1 0 1 1 1 0]'
$ws.Range("F3").Value = '[1 1 0 1 1 0 1 0 1 1 1 0]'
$ws.Range("G3").Value = 69
$ws.Range("H3").Value = 0.6624999999999999
$ws.Range("I3").Value = 0.02105927289454936
$ws.Range("J3").Value = 0.5540816326530613
$ws.Range("K3").Value = 0.0786473555337113

# Row 4
$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', None),
                (''selector'',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                (''model'',
                 BaggingClassifier(estimator=LogisticRegression(C=0.0001,
                                                                class_weight=''balanced'',
                                                                max_iter=1000,
                                                                random_state=42,
                                                                solver=''saga''),
                                   n_estimators=50, random_state=42))])'
$ws.Range("B4").Value = 0.6285714285714286
$ws.Range("C4").Value = '{''selector'': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), ''scaler'': None, ''model__n_estimators'': 50, ''model__estimator__solver'': ''saga'', ''model__estimator__penalty'': ''l2'', ''model__estimator__class_weight'': ''balanced'', ''model__estimator__C'': 0.0001}'
$ws.Range("D4").Value = 0.6666666666666666
$ws.Range("E4").Value = '[1 0 1 1 1 1 0 1 0 1 0 1]'
$ws.Range("F4").Value = '[1 0 0 1 1 1 1 0 1 0 0 1]'
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.6681385281385283
$ws.Range("I4").Value = 0.03433622428519705
$ws.Range("J4").Value = 0.5418181818181819
$ws.Range("K4").Value = 0.07881339753750574

# Row 5
$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', None),
                (''selector'',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'',
                                                     random_state=42))),
                (''model'',
                 BaggingClassifier(estimator=LogisticRegression(C=0.0001,
                                                                class_weight=''balanced'',
                                                                max_iter=1000,
                                                                random_state=42,
                                                                solver=''saga''),
                                   n_estimators=50, random_state=42))])'
$ws.Range("B5").Value = 0.6285714285714286
$ws.Range("C5").Value = '{''selector'': SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'', random_state=42)), ''scaler'': None, ''model__n_estimators'': 50, ''model__estimator__solver'': ''saga'', ''model__estimator__penalty'': ''l2'', ''model__estimator__class_weight'': ''balanced'', ''model__estimator__C'': 0.0001}'
$ws.Range("D5").Value = 0.75
$ws.Range("E5").Value = '[1 1 0 0 0 0 1 0 1 1 1 1]'
$ws.Range("F5").Value = '[1 1 1 0 0 1 1 1 1 0 1 1]'
$ws.Range("G5").Value = 11
$ws.Range("H5").Value = 0.6497795414462082
$ws.Range("I5").Value = 0.03880250498747913
$ws.Range("J5").Value = 0.5047619047619049
$ws.Range("K5").Value = 0.07708539671173281

# Row 6
$ws.Range("A6").Value = 'Pipeline(steps=[(''scaler'', None),
                (''selector'',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'',
                                                     random_state=42))),
                (''model'',
                 BaggingClassifier(estimator=LogisticRegression(C=1,
                                                                max_iter=1000,
                                                                penalty=''l1'',
                                                                random_state=42,
                                                                solver=''saga''),
                                   n_estimators=50, random_state=42))])'
$ws.Range("B6").Value = 0.6476190476190476
$ws.Range("C6").Value = '{''selector'': SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'', random_state=42)), ''scaler'': None, ''model__n_estimators'': 50, ''model__estimator__solver'': ''saga'', ''model__estimator__penalty'': ''l1'', ''model__estimator__class_weight'': None, ''model__estimator__C'': 1}'
$ws.Range("D6").Value = 0.5714285714285715
$ws.Range("E6").Value = '[1 1 1 1 0 0 0 0 1 1 0 0]'
$ws.Range("F6").Value = '[1 1 1 0 0 1 1 1 1 0 1 0]'
$ws.Range("G6").Value = 14
$ws.Range("H6").Value = 0.6643707482993195
$ws.Range("I6").Value = 0.0329769197326136
$ws.Range("J6").Value = 0.5627551020408164
$ws.Range("K6").Value = 0.07154926691588064
